$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.180165333333334
$ws.Range("H2").Value = 6.540496
$ws.Range("I2").Value = 0.01970539991828544
$ws.Range("J2").Value = 0.01970539991828544
$ws.Range("M2").Value = 1.910418
$ws.Range("N2").Value = 5.731254
$ws.Range("O2").Value = 0.01809124304049503
$ws.Range("P2").Value = 0.01809124304049503
$ws.Range("Q2").Value = 4.165027095776
$ws.Range("R2").Value = 37.485243861984
$ws.Range("S2").Value = 0.0003564951791318527
$ws.Range("T2").Value = 0.0003564951791318527
$ws.Range("G3").Value = 2.180165333333334
$ws.Range("H3").Value = 6.540496
$ws.Range("I3").Value = 0.01970539991828544
$ws.Range("J3").Value = 0.01970539991828544
$ws.Range("O3").Value = 0.302988173785169
$ws.Range("P3").Value = 0.302988173785169
$ws.Range("Q3").Value = 69.75496104331734
$ws.Range("R3").Value = 627.7946493898561
$ws.Range("S3").Value = 0.005970503134947723
$ws.Range("T3").Value = 0.005970503134947723
$ws.Range("G4").Value = 2.180165333333334
$ws.Range("H4").Value = 6.540496
$ws.Range("I4").Value = 0.01970539991828544
$ws.Range("J4").Value = 0.01970539991828544
$ws.Range("M4").Value = 37.858701
$ws.Range("N4").Value = 113.576103
$ws.Range("O4").Value = 0.3585136661130873
$ws.Range("P4").Value = 0.3585136661130873
$ws.Range("Q4").Value = 82.538227485232
$ws.Range("R4").Value = 742.8440473670879
$ws.Range("S4").Value = 0.007064655166929043
$ws.Range("T4").Value = 0.007064655166929043
$ws.Range("G5").Value = 2.180165333333334
$ws.Range("H5").Value = 6.540496
$ws.Range("I5").Value = 0.01970539991828544
$ws.Range("J5").Value = 0.01970539991828544
$ws.Range("M5").Value = 33.83466466666667
$ws.Range("N5").Value = 101.503994
$ws.Range("O5").Value = 0.3204069170612486
$ws.Range("P5").Value = 0.3204069170612486
$ws.Range("Q5").Value = 73.7651629712249
$ws.Range("R5").Value = 663.886466741024
$ws.Range("S5").Value = 0.006313746437276816
$ws.Range("T5").Value = 0.006313746437276816
$ws.Range("I6").Value = 0.733713204346044
$ws.Range("J6").Value = 0.7337132043460441
$ws.Range("M6").Value = 1.910418
$ws.Range("N6").Value = 5.731254
$ws.Range("O6").Value = 0.01809124304049503
$ws.Range("P6").Value = 0.01809124304049503
$ws.Range("Q6").Value = 155.081114278436
$ws.Range("R6").Value = 1395.730028505924
$ws.Range("S6").Value = 0.01327378390184467
$ws.Range("T6").Value = 0.01327378390184468
$ws.Range("I7").Value = 0.733713204346044
$ws.Range("J7").Value = 0.7337132043460441
$ws.Range("O7").Value = 0.302988173785169
$ws.Range("P7").Value = 0.302988173785169
$ws.Range("S7").Value = 0.2223064238668724
$ws.Range("T7").Value = 0.2223064238668724
$ws.Range("I8").Value = 0.733713204346044
$ws.Range("J8").Value = 0.7337132043460441
$ws.Range("M8").Value = 37.858701
$ws.Range("N8").Value = 113.576103
$ws.Range("O8").Value = 0.3585136661130873
$ws.Range("P8").Value = 0.3585136661130873
$ws.Range("Q8").Value = 3073.238179400602
$ws.Range("R8").Value = 27659.14361460541
$ws.Range("S8").Value = 0.263046210765681
$ws.Range("T8").Value = 0.2630462107656811
$ws.Range("I9").Value = 0.733713204346044
$ws.Range("J9").Value = 0.7337132043460441
$ws.Range("M9").Value = 33.83466466666667
$ws.Range("N9").Value = 101.503994
$ws.Range("O9").Value = 0.3204069170612486
$ws.Range("P9").Value = 0.3204069170612486
$ws.Range("Q9").Value = 2746.580851805152
$ws.Range("R9").Value = 24719.22766624637
$ws.Range("S9").Value = 0.2350867858116459
$ws.Range("T9").Value = 0.2350867858116459
$ws.Range("G10").Value = 25.672264
$ws.Range("H10").Value = 77.016792
$ws.Range("I10").Value = 0.2320384702908474
$ws.Range("J10").Value = 0.2320384702908474
$ws.Range("M10").Value = 1.910418
$ws.Range("N10").Value = 5.731254
$ws.Range("O10").Value = 0.01809124304049503
$ws.Range("P10").Value = 0.01809124304049503
$ws.Range("Q10").Value = 49.044755246352
$ws.Range("R10").Value = 441.402797217168
$ws.Range("S10").Value = 0.004197864360776405
$ws.Range("T10").Value = 0.004197864360776405
$ws.Range("G11").Value = 25.672264
$ws.Range("H11").Value = 77.016792
$ws.Range("I11").Value = 0.2320384702908474
$ws.Range("J11").Value = 0.2320384702908474
$ws.Range("O11").Value = 0.302988173785169
$ws.Range("P11").Value = 0.302988173785169
$ws.Range("Q11").Value = 821.3908128131679
$ws.Range("R11").Value = 7392.517315318511
$ws.Range("S11").Value = 0.07030491236132805
$ws.Range("T11").Value = 0.07030491236132805
$ws.Range("G12").Value = 25.672264
$ws.Range("H12").Value = 77.016792
$ws.Range("I12").Value = 0.2320384702908474
$ws.Range("J12").Value = 0.2320384702908474
$ws.Range("M12").Value = 37.858701
$ws.Range("N12").Value = 113.576103
$ws.Range("O12").Value = 0.3585136661130873
$ws.Range("P12").Value = 0.3585136661130873
$ws.Range("Q12").Value = 971.9185667690639
$ws.Range("R12").Value = 8747.267100921574
$ws.Range("S12").Value = 0.0831889626632444
$ws.Range("T12").Value = 0.0831889626632444
$ws.Range("G13").Value = 25.672264
$ws.Range("H13").Value = 77.016792
$ws.Range("I13").Value = 0.2320384702908474
$ws.Range("J13").Value = 0.2320384702908474
$ws.Range("M13").Value = 33.83466466666667
$ws.Range("N13").Value = 101.503994
$ws.Range("O13").Value = 0.3204069170612486
$ws.Range("P13").Value = 0.3204069170612486
$ws.Range("Q13").Value = 868.6124436741387
$ws.Range("R13").Value = 7817.511993067248
$ws.Range("S13").Value = 0.07434673090549854
$ws.Range("T13").Value = 0.07434673090549854
$ws.Range("G14").Value = 1.608999666666667
$ws.Range("H14").Value = 4.826999
$ws.Range("I14").Value = 0.01454292544482312
$ws.Range("J14").Value = 0.01454292544482312
$ws.Range("M14").Value = 1.910418
$ws.Range("N14").Value = 5.731254
$ws.Range("O14").Value = 0.01809124304049503
$ws.Range("P14").Value = 0.01809124304049503
$ws.Range("Q14").Value = 3.073861925194
$ws.Range("R14").Value = 27.664757326746
$ws.Range("S14").Value = 0.0002630995987420944
$ws.Range("T14").Value = 0.0002630995987420944
$ws.Range("G15").Value = 1.608999666666667
$ws.Range("H15").Value = 4.826999
$ws.Range("I15").Value = 0.01454292544482312
$ws.Range("J15").Value = 0.01454292544482312
$ws.Range("O15").Value = 0.302988173785169
$ws.Range("P15").Value = 0.302988173785169
$ws.Range("Q15").Value = 51.48036589291267
$ws.Range("R15").Value = 463.323293036214
$ws.Range("S15").Value = 0.004406334422020825
$ws.Range("T15").Value = 0.004406334422020825
$ws.Range("G16").Value = 1.608999666666667
$ws.Range("H16").Value = 4.826999
$ws.Range("I16").Value = 0.01454292544482312
$ws.Range("J16").Value = 0.01454292544482312
$ws.Range("M16").Value = 37.858701
$ws.Range("N16").Value = 113.576103
$ws.Range("O16").Value = 0.3585136661130873
$ws.Range("P16").Value = 0.3585136661130873
$ws.Range("Q16").Value = 60.91463728943299
$ws.Range("R16").Value = 548.2317356048969
$ws.Range("S16").Value = 0.00521383751723284
$ws.Range("T16").Value = 0.00521383751723284
$ws.Range("G17").Value = 1.608999666666667
$ws.Range("H17").Value = 4.826999
$ws.Range("I17").Value = 0.01454292544482312
$ws.Range("J17").Value = 0.01454292544482312
$ws.Range("M17").Value = 33.83466466666667
$ws.Range("N17").Value = 101.503994
$ws.Range("O17").Value = 0.3204069170612486
$ws.Range("P17").Value = 0.3204069170612486
$ws.Range("Q17").Value = 54.43996417044512
$ws.Range("R17").Value = 489.959677534006
$ws.Range("S17").Value = 0.004659653906827365
$ws.Range("T17").Value = 0.004659653906827365
